$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-27"

# Update the March row label text
$ws.Range("A4").Value = "March (through 03-27)"

# Update March row values (row 4)
$ws.Range("B4").Value = 24
$ws.Range("C4").Value = 38
$ws.Range("D4").Value = 49
$ws.Range("E4").Value = 55
$ws.Range("H4").Value = 75
$ws.Range("I4").Value = 111

# Update Total row values (row 5)
$ws.Range("B5").Value = 61
$ws.Range("C5").Value = 125
$ws.Range("D5").Value = 180
$ws.Range("E5").Value = 192
$ws.Range("H5").Value = 417
$ws.Range("I5").Value = 411
